# Adding area to the discharge file (Station3_2021-06-18_1300.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header labels: Area (col G) and Atotal (col H, and mirrored in J/K) ---
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- Fill in the last missing segment's raw measurements (row 9) ---
$ws.Range("A9").Value = 40
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0

# --- Area per segment, mirroring the existing Q (discharge) column logic ---
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Totals ---
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# --- Summary cells pulling the two totals together ---
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Update the active selection to match the edited workbook ---
$ws.Range("J2:K2").Select()

$wb.Save()
